# Generate Report for Handoff
# The "1d6615c0-25c9-4c37-b0d4-5c4cbaf1d179.md" file (row 3 on every sheet)
# has finished translation and is now ready to be handed off. Update its
# status + handoff timestamps on the Overview sheet and on each of the
# per-language detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-13-14 08:13:33"

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-14 08:13:30"

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-14 08:13:33"
